# Update Name of Algo
# Apply updated KNN imputation results to column C for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = -12.342
$ws.Range("C10").Value = -12.179
$ws.Range("C12").Value = -12.157
$ws.Range("C18").Value = -12.157
$ws.Range("C25").Value = -12.324
